$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'302.64"
$ws.Range("E2").Value = "'0.76%"

$ws.Range("D3").Value = "'32.16"
$ws.Range("E3").Value = "'1.18%"

$ws.Range("D4").Value = "'4.998"
$ws.Range("E4").Value = "'-2.86%"

$ws.Range("D5").Value = "'0.07914"
$ws.Range("E5").Value = "'-2.78%"

$ws.Range("D6").Value = "'2.111"
$ws.Range("E6").Value = "'-16.15%"

$ws.Range("D7").Value = "'7.868"
$ws.Range("E7").Value = "'0.29%"

$ws.Range("D8").Value = "'3.804"
$ws.Range("E8").Value = "'-1.73%"

$ws.Range("D9").Value = "'0.9261"
$ws.Range("E9").Value = "'0.07%"

$ws.Range("E10").Value = "'-0.16%"

$ws.Range("D11").Value = "'0.07958"
$ws.Range("E11").Value = "'7.12%"

$ws.Range("D12").Value = "'0.08760"
$ws.Range("E12").Value = "'-1.46%"

$ws.Range("D13").Value = "'0.03168"
$ws.Range("E13").Value = "'4.89%"

$ws.Range("E14").Value = "'0.28%"

$ws.Range("D15").Value = "'0.001518"
$ws.Range("E15").Value = "'-0.72%"

$ws.Range("D16").Value = "'0.006024"
$ws.Range("E16").Value = "'0.31%"

$ws.Range("D17").Value = "'3.467"
$ws.Range("E17").Value = "'-3.90%"

$ws.Range("D18").Value = "'2.278"
$ws.Range("E18").Value = "'-0.30%"

$ws.Range("E19").Value = "'0.77%"

$ws.Range("D20").Value = "'0.1291"
$ws.Range("E20").Value = "'-4.38%"

$ws.Range("D21").Value = "'4.198"
$ws.Range("E21").Value = "'-0.85%"

$ws.Range("D22").Value = "'0.1792"
$ws.Range("E22").Value = "'6.72%"

$ws.Range("D23").Value = "'0.04611"
$ws.Range("E23").Value = "'-0.30%"

$ws.Range("D24").Value = "'0.001239"
$ws.Range("E24").Value = "'-0.56%"

$ws.Range("D25").Value = "'0.004484"
$ws.Range("E25").Value = "'-1.36%"

$ws.Range("D26").Value = "'0.0001251"
$ws.Range("E26").Value = "'4.41%"

$ws.Range("D39").Value = "'0.01734"
$ws.Range("E39").Value = "'-2.28%"

$ws.Range("D40").Value = "'0.04817"
$ws.Range("E40").Value = "'4.42%"

$ws.Range("D41").Value = "'0.007335"
$ws.Range("E41").Value = "'5.99%"

$ws.Range("E42").Value = "'-0.70%"

$ws.Range("D43").Value = "'0.002362"
$ws.Range("E43").Value = "'10.55%"

$ws.Range("D44").Value = "'0.01110"
$ws.Range("E44").Value = "'12.69%"

$ws.Range("D45").Value = "'0.00006026"
$ws.Range("E45").Value = "'-2.56%"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.22%"

$ws.Range("D47").Value = "'0.003393"
$ws.Range("E47").Value = "'-59.55%"

$ws.Range("D48").Value = "'0.8204"
$ws.Range("E48").Value = "'2.22%"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("E49").Value = "'0.22%"

$ws.Range("D50").Value = "'0.0002002"
$ws.Range("E50").Value = "'0.22%"
